$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$targetRows = @(3,6,7,9,15,19,20,23,29,30,34,36,39,46,51,56,58,64,69,72,73,75,77,79,83,86,87,89,91,93,95,99,103,105,107,108,111,112,116,123,124,125,126,132,134,137,142,145,146,150,153,156,157,159,164,166,173,177,180,183,185,186,188,193)

foreach ($r in $targetRows) {
    $ws.Range("G$r`:N$r").Value = "N/A"
}
